# Dataschema_TRACY.xlsx - "Categories" sheet
# Add a "don't know" category (value 2) to the two family-history
# variables FAM1_CHD_STROKE and FAM1_CANCER, inserted right after each
# variable's existing "no"/"yes" rows. Everything below shifts down by
# one row for each insertion (two rows total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# --- FAM1_CHD_STROKE -------------------------------------------------
# Existing rows 39 (no/0) and 40 (yes/1) stay put; insert a fresh row
# at 41 for the new "don't know" / 2 category.
$ws.Rows.Item(41).Insert()
$ws.Cells.Item(41, 1).Value = "FAM1_CHD_STROKE"
$ws.Cells.Item(41, 2).Value = "don't know"
$ws.Cells.Item(41, 3).Value = 2

# --- FAM1_CANCER -------------------------------------------------------
# After the insert above, FAM1_CANCER's "no"/"yes" rows now sit at 44
# and 45. Insert a fresh row at 46 for its "don't know" / 2 category.
$ws.Rows.Item(46).Insert()
$ws.Cells.Item(46, 1).Value = "FAM1_CANCER"
$ws.Cells.Item(46, 2).Value = "don't know"
$ws.Cells.Item(46, 3).Value = 2
